# Update weekly fruit/vegetable price records (Haba - Terminal La Palmera de La Serena)
# Rows 2 and 4-14 get their Fecha/Volumen/Precio.../Origen values re-shuffled to reflect
# the latest weekly data pull. Row 3 is untouched.

$ws = $excel.ActiveWorkbook.ActiveSheet

$rows = @{
    2  = @{ D = 44446; J = 500; K = 11000; L = 12000; M = 11500; O = "Provincia del Elquí"; P = 460 }
    4  = @{ D = 44425; J = 400; K = 11500; L = 12000; M = 11750; O = "Provincia del Elquí"; P = 470 }
    5  = @{ D = 44356; J = 500; K = 13000; L = 14000; M = 13500; O = "Provincia de Limarí";  P = 540 }
    6  = @{ D = 44376; J = 400; K = 12000; L = 13000; M = 12500; O = "Provincia del Elquí"; P = 500 }
    7  = @{ D = 44484; J = 400; K = 9000;  L = 10000; M = 9500;  O = "Provincia del Elquí"; P = 380 }
    8  = @{ D = 44384; J = 560; K = 11500; L = 12000; M = 11750; O = "Provincia del Elquí"; P = 470 }
    9  = @{ D = 44370; J = 520; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
    10 = @{ D = 44386; J = 500; K = 11000; L = 12000; M = 11500; O = "Provincia del Elquí"; P = 460 }
    11 = @{ D = 44473; J = 500; K = 8500;  L = 9000;  M = 8750;  O = "Provincia del Elquí"; P = 350 }
    12 = @{ D = 44377; J = 520; K = 12500; L = 13000; M = 12750; O = "Provincia del Elquí"; P = 510 }
    13 = @{ D = 44316; J = 300; K = 16000; L = 17000; M = 16500; O = "Provincia del Elquí"; P = 660 }
    14 = @{ D = 44372; J = 500; K = 13000; L = 14000; M = 13500; O = "Provincia del Elquí"; P = 540 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Cells.Item($r, 4).Value2  = $vals.D   # D - Fecha
    $ws.Cells.Item($r, 10).Value2 = $vals.J   # J - Volumen
    $ws.Cells.Item($r, 11).Value2 = $vals.K   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value2 = $vals.L   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value2 = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 15).Value  = $vals.O   # O - Origen
    $ws.Cells.Item($r, 16).Value2 = $vals.P   # P - Precio $/Kg
}
